# Apply cryptos list update (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D that need to stay as text even though the new value looks numeric:
# force a text number format first so Excel does not auto-convert them to numbers.
$forceTextCells = @("D5", "D6", "D8", "D9", "D11", "D15", "D18", "D19", "D22", "D26", "D28", "D31", "D32", "D33", "D36", "D38", "D41", "D46", "D47", "D48", "D49")
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "60.918.40"

# Row 3
$ws.Range("D3").Value = "2.637.16"
$ws.Range("E3").Value = "  +1.68%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").Value = "529.35"
$ws.Range("E5").Value = "  +3.97%  "

# Row 6
$ws.Range("D6").Value = "155.08"
$ws.Range("E6").Value = "  +0.95%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").Value = "0.590"
$ws.Range("E8").Value = "  +0.40%  "

# Row 9
$ws.Range("D9").Value = "6.64"
$ws.Range("E9").Value = "  -0.62%  "

# Row 10
$ws.Range("E10").Value = "  +5.13%  "

# Row 11
$ws.Range("D11").Value = "0.352"
$ws.Range("E11").Value = "  +1.84%  "

# Row 12
$ws.Range("E12").Value = "  -0.08%  "

# Row 13
$ws.Range("D13").Value = "3.098.87"
$ws.Range("E13").Value = "  +1.63%  "

# Row 14
$ws.Range("D14").Value = "60.932.52"
$ws.Range("E14").Value = "  +1.04%  "

# Row 15
$ws.Range("D15").Value = "21.91"
$ws.Range("E15").Value = "  +2.04%  "

# Row 16
$ws.Range("E16").Value = "  +2.69%  "

# Row 17
$ws.Range("D17").Value = "2.636.27"
$ws.Range("E17").Value = "  +1.45%  "

# Row 18
$ws.Range("D18").Value = "4.75"
$ws.Range("E18").Value = "  +0.42%  "

# Row 19
$ws.Range("D19").Value = "353.69"
$ws.Range("E19").Value = "  +0.17%  "

# Row 20
$ws.Range("E20").Value = "  +1.27%  "

# Row 21
$ws.Range("E21").Value = "  +1.63%  "

# Row 22
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.16%  "

# Row 23
$ws.Range("E23").Value = "  +2.00%  "

# Row 24
$ws.Range("E24").Value = "  +2.44%  "

# Row 25
$ws.Range("E25").Value = "  +1.50%  "

# Row 26
$ws.Range("D26").Value = "0.991"
$ws.Range("E26").Value = "  -0.65%  "

# Row 27
$ws.Range("D27").Value = "0.0₃0863"
$ws.Range("E27").Value = "  +3.39%  "

# Row 28
$ws.Range("D28").Value = "7.39"
$ws.Range("E28").Value = "  +0.60%  "

# Row 29
$ws.Range("E29").Value = "  -0.01%  "

# Row 30
$ws.Range("E30").Value = "  +7.43%  "

# Row 31
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "1.63"
$ws.Range("E31").Value = "  +4.31%  "

# Row 32
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "19.46"
$ws.Range("E32").Value = "  +0.67%  "

# Row 33
$ws.Range("D33").Value = "150.31"
$ws.Range("E33").Value = "  -0.57%  "

# Row 34
$ws.Range("E34").Value = "  +4.10%  "

# Row 35
$ws.Range("E35").Value = "  +1.86%  "

# Row 36
$ws.Range("D36").Value = "0.925"
$ws.Range("E36").Value = "  +10.44%  "

# Row 37
$ws.Range("E37").Value = "  +2.28%  "

# Row 38
$ws.Range("D38").Value = "307.32"
$ws.Range("E38").Value = "  +4.31%  "

# Row 39
$ws.Range("E39").Value = "  +1.58%  "

# Row 40
$ws.Range("E40").Value = "  +1.84%  "

# Row 41
$ws.Range("D41").Value = "36.58"
$ws.Range("E41").Value = "  +0.99%  "

# Row 42
$ws.Range("E42").Value = "  +4.00%  "

# Row 43
$ws.Range("E43").Value = "  +1.74%  "

# Row 44
$ws.Range("E44").Value = "  +1.88%  "

# Row 45
$ws.Range("E45").Value = "  -0.04%  "

# Row 46
$ws.Range("D46").Value = "19.78"
$ws.Range("E46").Value = "  +0.81%  "

# Row 47
$ws.Range("D47").Value = "4.95"
$ws.Range("E47").Value = "  +3.21%  "

# Row 48
$ws.Range("D48").Value = "0.0239"
$ws.Range("E48").Value = "  +2.44%  "

# Row 49
$ws.Range("D49").Value = "19.31"
$ws.Range("E49").Value = "  +8.35%  "

# Row 50
$ws.Range("E50").Value = "  +0.29%  "

# Row 51
$ws.Range("D51").Value = "1.981.64"
$ws.Range("E51").Value = "  -0.30%  "

# Restore default (General) style for the forced-text cells so formatting matches the source data
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).Style = "Normal"
}
